$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 96
$ws1.Range("F4").Value = 1205
$ws1.Range("F5").Value = 607

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 96
$ws4.Range("F4").Value = 1205
$ws4.Range("F6").Value = 607
